# Stricter separation of departement specific pre mid and post mid courses
#
# Updates the Classroom_Allocation sheet so rooms allocated to "post mid"
# electives (ELECTIVE_B8 / ELECTIVE_B9 baskets) no longer overlap with the
# rooms used for the corresponding "pre mid" baskets (ELECTIVE_B6 / B7),
# and refreshes the two summary sheets (Basket_Course_Allocations,
# Executive_Summary) that mirror that allocation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Classroom_Allocation sheet: room / capacity / facilities / room-code
# reassignments for the affected rows.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Allocation")

# Row 5
$ws.Range("G5").Value = "Auditorium"
$ws.Range("H5").Value = "240"
$ws.Range("I5").Value = "Audio/Video System"
$ws.Range("M5").Value = "C004"

# Row 17
$ws.Range("G17").Value = "Auditorium"
$ws.Range("H17").Value = "240"
$ws.Range("I17").Value = "Audio/Video System"
$ws.Range("M17").Value = "C004"

# Row 18
$ws.Range("G18").Value = "large classroom"
$ws.Range("H18").Value = "120"
$ws.Range("I18").Value = ""
$ws.Range("M18").Value = "C001"

# Row 19
$ws.Range("G19").Value = "large classroom"
$ws.Range("H19").Value = "120"
$ws.Range("M19").Value = "C002"

# Row 20
$ws.Range("M20").Value = "C101"

# Row 21
$ws.Range("G21").Value = "Auditorium"
$ws.Range("H21").Value = "240"
$ws.Range("I21").Value = "Audio/Video System"
$ws.Range("M21").Value = "C004"

# Row 37
$ws.Range("G37").Value = "Auditorium"
$ws.Range("H37").Value = "240"
$ws.Range("I37").Value = "Audio/Video System"
$ws.Range("M37").Value = "C004"

# Row 38
$ws.Range("G38").Value = "large classroom"
$ws.Range("H38").Value = "120"
$ws.Range("M38").Value = "C002"

# Row 39
$ws.Range("M39").Value = "C102"

# Row 40
$ws.Range("I40").Value = "Projector"
$ws.Range("M40").Value = "C104"

# Row 41
$ws.Range("M41").Value = "C203"

# Row 42
$ws.Range("G42").Value = "Auditorium"
$ws.Range("H42").Value = "240"
$ws.Range("I42").Value = "Audio/Video System"
$ws.Range("M42").Value = "C004"

# Row 43
$ws.Range("G43").Value = "large classroom"
$ws.Range("H43").Value = "120"
$ws.Range("I43").Value = ""
$ws.Range("M43").Value = "C001"

# Row 44
$ws.Range("G44").Value = "large classroom"
$ws.Range("H44").Value = "120"
$ws.Range("M44").Value = "C002"

# Row 45
$ws.Range("M45").Value = "C101"

# Row 46
$ws.Range("G46").Value = "Auditorium"
$ws.Range("H46").Value = "240"
$ws.Range("I46").Value = "Audio/Video System"
$ws.Range("M46").Value = "C004"

# Row 47
$ws.Range("G47").Value = "large classroom"
$ws.Range("H47").Value = "120"
$ws.Range("I47").Value = ""
$ws.Range("M47").Value = "C001"

# Row 48
$ws.Range("G48").Value = "large classroom"
$ws.Range("H48").Value = "120"
$ws.Range("M48").Value = "C002"

# Row 49
$ws.Range("I49").Value = "Projector"
$ws.Range("M49").Value = "C101"

# ---------------------------------------------------------------------
# Basket_Course_Allocations sheet: each basket row now lists the single
# room actually reserved (instead of the old "both candidate rooms"
# list), matching the tighter per-basket room assignment above.
# ---------------------------------------------------------------------
$wsBasket = $wb.Worksheets.Item("Basket_Course_Allocations")

$wsBasket.Range("C9").Value = "C004"
$wsBasket.Range("C10").Value = "C002"
$wsBasket.Range("C11").Value = "C102"
$wsBasket.Range("C12").Value = "C104"
$wsBasket.Range("C13").Value = "C203"
$wsBasket.Range("C14").Value = "C004"
$wsBasket.Range("C15").Value = "C001"
$wsBasket.Range("C16").Value = "C002"
$wsBasket.Range("C17").Value = "C101"

# ---------------------------------------------------------------------
# Executive_Summary sheet: regeneration timestamp bump.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Executive_Summary")
$wsSummary.Range("C3").Value = "2026-01-26 12:46"
